$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 305
$ws.Range("I2").Value = 850
$ws.Range("J2").Value = 3495
$ws.Range("K2").Value = 16
$ws.Range("L2").Value = 930
$ws.Range("M2").Value = 51
$ws.Range("N2").Value = 602
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 10
$ws.Range("Q2").Value = 7
$ws.Range("R2").Value = 48
$ws.Range("S2").Value = 414
$ws.Range("T2").Value = 669
$ws.Range("U2").Value = 38
$ws.Range("V2").Value = 5466
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 5430
$ws.Range("Y2").Value = 6
$ws.Range("Z2").Value = 77
$ws.Range("AA2").Value = 41
